$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 105, shifting existing rows 105:136 down to 106:137
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with data (mirrors the constant columns used throughout this block)
$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44785
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100102
$ws.Cells.Item(105, 8).Value = "Cítricos"
$ws.Cells.Item(105, 9).Value = 100102004
$ws.Cells.Item(105, 10).Value = "Mandarina"
$ws.Cells.Item(105, 11).Value = "Murcott"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 270
$ws.Cells.Item(105, 14).Value = 7000
$ws.Cells.Item(105, 15).Value = 7500
$ws.Cells.Item(105, 16).Value = 7278
$ws.Cells.Item(105, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(105, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(105, 19).Value = 404
$ws.Cells.Item(105, 20).Value = 18
